$wb = $excel.ActiveWorkbook

# Map: sheet name -> { row -> new F value }
$updates = @{
    "展览" = @{
        3  = 153
        4  = 68
        5  = 510
        6  = 1522
        7  = 1065
        8  = 112
        9  = 214
        10 = 154
        11 = 220
        12 = 121
        13 = 184
        14 = 169
    }
    "全部类型" = @{
        3  = 153
        4  = 68
        5  = 510
        6  = 1522
        8  = 1065
        9  = 112
        10 = 214
        11 = 154
        12 = 220
        13 = 121
        14 = 184
        15 = 169
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rowsMap = $updates[$sheetName]
    foreach ($row in $rowsMap.Keys) {
        $ws.Cells.Item($row, 6).Value = $rowsMap[$row]
    }
}
